$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text values are kept as text (matching original formatting),
# since these price/volume columns store data as text, not numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.619.73"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.61%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.228.06"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.02%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.24%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "270.36"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +4.58%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "89.63"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +11.75%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.621"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.57%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.13%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.606"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.23%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "45.72"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +6.32%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0918"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.56%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.88"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +12.18%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.10%  "

$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.09"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.09%  "

$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.559.12"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.44%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.206.19"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.40%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.791"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.61%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.576.61"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.54%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.87%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.33"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.06%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.97"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.91%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.34"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.04%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "231.81"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.14%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.62"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -8.90%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.01%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.48"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +11.62%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.90"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.05%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +5.39%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.83%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.62"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -4.86%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "172.38"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.10%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0907"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.83%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.70"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.91%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.35"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.01%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.24%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -4.36%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0350"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -4.75%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.26"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -5.30%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +12.88%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.48"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.71%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.15"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.96%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.212"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +5.06%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "63.06"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.19%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.33"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.08%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.51"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.48%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.19%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "100.04"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.57%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.18%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.15%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.435"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.59%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.49"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.07%  "
